$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / shared-string edits (header banner + week-range) ---
$ws.Range("A8").Characters(21, 2).Text = "29"
$ws.Range("C9").Characters(27, 8).Text = "7/15/2024"
$ws.Range("C9").Characters(47, 9).Text = "7/21/2024"

# --- Numeric grid edits (weekly crime-stat refresh) ---
$ws.Range("M14").Value = -38.461538461538
$ws.Range("N14").Value = -63.636363636363
$ws.Range("D15").Value = 1
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 25
$ws.Range("K15").Value = -28
$ws.Range("M15").Value = -10
$ws.Range("N15").Value = -68.421052631578
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 7
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 157
$ws.Range("J16").Value = 144
$ws.Range("K16").Value = 9.027777777777
$ws.Range("L16").Value = -19.897959183673
$ws.Range("M16").Value = -31.140350877193
$ws.Range("N16").Value = -88.327137546468
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 22
$ws.Range("E17").Value = -68.181818181818
$ws.Range("F17").Value = 71
$ws.Range("G17").Value = 61
$ws.Range("H17").Value = 16.393442622950
$ws.Range("I17").Value = 450
$ws.Range("J17").Value = 379
$ws.Range("K17").Value = 18.733509234828
$ws.Range("L17").Value = 25
$ws.Range("M17").Value = 79.282868525896
$ws.Range("N17").Value = -35.530085959885
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -25
$ws.Range("I18").Value = 116
$ws.Range("J18").Value = 114
$ws.Range("K18").Value = 1.754385964912
$ws.Range("L18").Value = -0.854700854700
$ws.Range("M18").Value = -43.137254901960
$ws.Range("N18").Value = -89.868995633187
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -23.076923076923
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -2.127659574468
$ws.Range("I19").Value = 330
$ws.Range("J19").Value = 383
$ws.Range("K19").Value = -13.838120104438
$ws.Range("L19").Value = -17.910447761194
$ws.Range("M19").Value = 14.186851211072
$ws.Range("N19").Value = -22.352941176470
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 23
$ws.Range("H20").Value = -34.782608695652
$ws.Range("I20").Value = 117
$ws.Range("J20").Value = 142
$ws.Range("K20").Value = -17.605633802816
$ws.Range("L20").Value = -29.090909090909
$ws.Range("M20").Value = -26.415094339622
$ws.Range("N20").Value = -89.878892733564
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 57
$ws.Range("E21").Value = -43.859649122807
$ws.Range("F21").Value = 162
$ws.Range("G21").Value = 162
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 1196
$ws.Range("J21").Value = 1195
$ws.Range("K21").Value = 0.083682008368
$ws.Range("L21").Value = -5.454545454545
$ws.Range("M21").Value = 2.749140893470
$ws.Range("N21").Value = -75.330033003300
$ws.Range("M22").Value = 100
$ws.Range("G23").Value = 1
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = 21.739130434782
$ws.Range("F24").Value = 93
$ws.Range("G24").Value = 76
$ws.Range("H24").Value = 22.368421052631
$ws.Range("I24").Value = 704
$ws.Range("J24").Value = 657
$ws.Range("K24").Value = 7.153729071537
$ws.Range("L24").Value = 29.650092081031
$ws.Range("M24").Value = 38.310412573673
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 200
$ws.Range("F25").Value = 22
$ws.Range("G25").Value = 10
$ws.Range("H25").Value = 120
$ws.Range("I25").Value = 108
$ws.Range("J25").Value = 121
$ws.Range("K25").Value = -10.743801652892
$ws.Range("L25").Value = -3.571428571428
$ws.Range("C26").Value = 17
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = 21.428571428571
$ws.Range("F26").Value = 87
$ws.Range("G26").Value = 65
$ws.Range("H26").Value = 33.846153846153
$ws.Range("I26").Value = 515
$ws.Range("J26").Value = 424
$ws.Range("K26").Value = 21.462264150943
$ws.Range("L26").Value = 22.327790973871
$ws.Range("M26").Value = -1.151631477927
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 32
$ws.Range("K27").Value = -9.375
$ws.Range("C28").Value = 2
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -16.666666666666
$ws.Range("I28").Value = 44
$ws.Range("J28").Value = 41
$ws.Range("K28").Value = 7.317073170731
$ws.Range("L28").Value = 2.325581395348
$ws.Range("C29").Value = "'20"
$ws.Range("C14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("D29").Value = 2
$ws.Range("I14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 200
$ws.Range("J29").Value = 25
$ws.Range("K29").Value = -8
$ws.Range("L29").Value = -17.857142857142
$ws.Range("M29").Value = -41.025641025641
$ws.Range("N29").Value = -80.833333333333
$ws.Range("C30").Value = "'20"
$ws.Range("C14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("D30").Value = 2
$ws.Range("I14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 50
$ws.Range("J30").Value = 21
$ws.Range("K30").Value = -9.523809523809
$ws.Range("L30").Value = -13.636363636363
$ws.Range("M30").Value = -42.424242424242
$ws.Range("N30").Value = -83.333333333333
$ws.Range("D31").Value = "'20"
$ws.Range("C14").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = "'21"
$ws.Range("C14").Copy()
$ws.Range("E31").PasteSpecial(-4122)

$excel.CutCopyMode = 0

